$wb = $excel.ActiveWorkbook

# --- parametrosInicio: update the migration label date (B10) ---
$wsParam = $wb.Worksheets.Item("parametrosInicio")
$wsParam.Range("B10").Value = "MIGRACIONES SGV FEBRERO 2023 28.02.2023"

# --- Give the "Fecha" input cell (B5) a short-date number format ---
$wsParam.Range("B5").NumberFormat = "mm-dd-yy"

# --- Move the visible selection from C12 to B6 ---
[void]$wsParam.Range("B6").Select()
